$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 190. This shifts the
# existing rows 190-200 down to 192-202, matching the target layout where
# two brand-new observations are prepended to this block and everything
# else keeps its relative order.
$ws.Rows.Item(190).Insert()
$ws.Rows.Item(190).Insert()

# Common (unchanged) attributes shared by every row in this Pera /
# Packham's Triumph / Terminal Hortofruticola Agro Chillan block.
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$tipo = "Fruta"
$productoId = 100104
$producto = "Frutos de pepita"
$categoriaId = 100104005
$categoria = "Pera"
$variedad = "Packham's Triumph"
$unidad = "`$/caja 16 kilos empedrada"
$origen = "Provincia de Curicó"
$kgUnidad = 16

# New row 190: Fecha 44753, Calidad Especial, Volumen 50, precios 9000/9000/9000, 562 $/Kg
$ws.Cells.Item(190, 1).Value = $mercadoId
$ws.Cells.Item(190, 2).Value = $mercado
$ws.Cells.Item(190, 3).Value = $region
$ws.Cells.Item(190, 4).Value = 44753
$ws.Cells.Item(190, 5).Value = $codreg
$ws.Cells.Item(190, 6).Value = $tipo
$ws.Cells.Item(190, 7).Value = $productoId
$ws.Cells.Item(190, 8).Value = $producto
$ws.Cells.Item(190, 9).Value = $categoriaId
$ws.Cells.Item(190, 10).Value = $categoria
$ws.Cells.Item(190, 11).Value = $variedad
$ws.Cells.Item(190, 12).Value = "Especial"
$ws.Cells.Item(190, 13).Value = 50
$ws.Cells.Item(190, 14).Value = 9000
$ws.Cells.Item(190, 15).Value = 9000
$ws.Cells.Item(190, 16).Value = 9000
$ws.Cells.Item(190, 17).Value = $unidad
$ws.Cells.Item(190, 18).Value = $origen
$ws.Cells.Item(190, 19).Value = 562
$ws.Cells.Item(190, 20).Value = $kgUnidad

# New row 191: Fecha 44753, Calidad Primera, Volumen 80, precios 7500/8000/7750, 484 $/Kg
$ws.Cells.Item(191, 1).Value = $mercadoId
$ws.Cells.Item(191, 2).Value = $mercado
$ws.Cells.Item(191, 3).Value = $region
$ws.Cells.Item(191, 4).Value = 44753
$ws.Cells.Item(191, 5).Value = $codreg
$ws.Cells.Item(191, 6).Value = $tipo
$ws.Cells.Item(191, 7).Value = $productoId
$ws.Cells.Item(191, 8).Value = $producto
$ws.Cells.Item(191, 9).Value = $categoriaId
$ws.Cells.Item(191, 10).Value = $categoria
$ws.Cells.Item(191, 11).Value = $variedad
$ws.Cells.Item(191, 12).Value = "Primera"
$ws.Cells.Item(191, 13).Value = 80
$ws.Cells.Item(191, 14).Value = 7500
$ws.Cells.Item(191, 15).Value = 8000
$ws.Cells.Item(191, 16).Value = 7750
$ws.Cells.Item(191, 17).Value = $unidad
$ws.Cells.Item(191, 18).Value = $origen
$ws.Cells.Item(191, 19).Value = 484
$ws.Cells.Item(191, 20).Value = $kgUnidad
